$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correccion de error en get-data: actualiza columna Nacionalidad (D)
$rowsToTrue = @(4,8,9,11,13,15,16,17,19,23,30,32,33,35,36,37,38,39,40,41,42,43,44,45,46,47,49,50,51,52,53,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,81,82,85,87,89,94,95,98,99,101,103,104,105,106,107,108,110,112,113,114,116,117,118,119,120,123,124,125,128,129,132,133,135,136,138,140,141,142,143,144,145,146,149,150,151,152,153,154,156,158,160,163,164,167,168,169,171,172,173,174,178,179,181,184,185,189,192,193,194,195,196,197,198,199,200,201,202,203,204,205,206,207,209,211,212,213,215,216,217,219,221,222,223,224,225,226,228,229,230,232,233,234,235,236,237,239,240,241)
foreach ($r in $rowsToTrue) {
    $ws.Cells.Item($r, 4).Value = $true
}

$rowsToFalse = @(97,177,191)
foreach ($r in $rowsToFalse) {
    $ws.Cells.Item($r, 4).Value = $false
}
